# Applies the Zephyr test-case workbook edits described in the commit.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: "Verify amount filter and sort functionality" -> split into a new
# single-amount-filter test case with 4 steps (rows 2-5).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("B2").Value = "Verify single amount filter functionality"
$ws1.Range("D2").Value = "Click on amount filter option"
$ws1.Range("F2").Value = "Amount filter dropdown is displayed"
$ws1.Range("I2").Value = "User is logged in as investigator and on transaction list page`nAmount filter is successfully applied"

$ws1.Range("B3").Value = "Verify single amount filter functionality"
$ws1.Range("D3").Value = "Select 'Greater than $1000' option"
$ws1.Range("F3").Value = "Only transactions over $1000 are shown in the list"
$ws1.Range("I3").Value = "User is logged in as investigator and on transaction list page`nAmount filter is successfully applied"

$ws1.Cells.Item(4, 1).Value = 1
$ws1.Cells.Item(4, 2).Value = "Verify single amount filter functionality"
$ws1.Cells.Item(4, 3).Value = 3
$ws1.Cells.Item(4, 4).Value = "Check transaction count"
$ws1.Cells.Item(4, 5).Value = ""
$ws1.Cells.Item(4, 6).Value = "Count reflects number of filtered transactions"
$ws1.Cells.Item(4, 7).Value = "6414a0cd67102fc717c034d7"
$ws1.Cells.Item(4, 8).Value = "This test case has been built by GenAI Workbench for XL import via Internal Importer."
$ws1.Cells.Item(4, 9).Value = "User is logged in as investigator and on transaction list page`nAmount filter is successfully applied"
$ws1.Cells.Item(4, 10).Value = "Core"
$ws1.Cells.Item(4, 11).Value = "external"
$ws1.Cells.Item(4, 12).Value = "INVHUB-10821"
$ws1.Cells.Item(4, 13).Value = "blocks"
$ws1.Cells.Item(4, 14).Value = "GenAI_Test_Case"
$ws1.Cells.Item(4, 15).Value = "IM-5000"
$ws1.Cells.Item(4, 16).Value = "blocks"
$ws1.Cells.Item(4, 17).Value = "IM-3000"
$ws1.Cells.Item(4, 18).Value = "3 - Medium"
$ws1.Cells.Item(4, 19).Value = 37
$ws1.Cells.Item(4, 20).Value = "Release-1.0"
$ws1.Cells.Item(4, 21).Value = "Dublin"

$ws1.Cells.Item(5, 1).Value = 1
$ws1.Cells.Item(5, 2).Value = "Verify single amount filter functionality"
$ws1.Cells.Item(5, 3).Value = 4
$ws1.Cells.Item(5, 4).Value = "Verify filter indicator"
$ws1.Cells.Item(5, 5).Value = ""
$ws1.Cells.Item(5, 6).Value = "Filter indicator is visible in the UI"
$ws1.Cells.Item(5, 7).Value = "6414a0cd67102fc717c034d7"
$ws1.Cells.Item(5, 8).Value = "This test case has been built by GenAI Workbench for XL import via Internal Importer."
$ws1.Cells.Item(5, 9).Value = "User is logged in as investigator and on transaction list page`nAmount filter is successfully applied"
$ws1.Cells.Item(5, 10).Value = "Core"
$ws1.Cells.Item(5, 11).Value = "external"
$ws1.Cells.Item(5, 12).Value = "INVHUB-10821"
$ws1.Cells.Item(5, 13).Value = "blocks"
$ws1.Cells.Item(5, 14).Value = "GenAI_Test_Case"
$ws1.Cells.Item(5, 15).Value = "IM-5000"
$ws1.Cells.Item(5, 16).Value = "blocks"
$ws1.Cells.Item(5, 17).Value = "IM-3000"
$ws1.Cells.Item(5, 18).Value = "3 - Medium"
$ws1.Cells.Item(5, 19).Value = 37
$ws1.Cells.Item(5, 20).Value = "Release-1.0"
$ws1.Cells.Item(5, 21).Value = "Dublin"

# ---------------------------------------------------------------------------
# Sheet2: "Verify filter clearing on Copilot refresh" -> new multiple
# concurrent filters test case with 3 steps (rows 2-4).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("B2").Value = "Verify multiple concurrent filters functionality"
$ws2.Range("D2").Value = "Apply amount filter for 'Greater than $500'"
$ws2.Range("F2").Value = "Transactions over $500 are displayed"
$ws2.Range("I2").Value = "User is logged in as investigator and on transaction list page`nMultiple filters are successfully applied"

$ws2.Range("B3").Value = "Verify multiple concurrent filters functionality"
$ws2.Range("D3").Value = "Apply date filter for 'Last 7 days'"
$ws2.Range("F3").Value = "Only transactions matching both amount and date criteria are shown"
$ws2.Range("I3").Value = "User is logged in as investigator and on transaction list page`nMultiple filters are successfully applied"

$ws2.Cells.Item(4, 1).Value = 2
$ws2.Cells.Item(4, 2).Value = "Verify multiple concurrent filters functionality"
$ws2.Cells.Item(4, 3).Value = 3
$ws2.Cells.Item(4, 4).Value = "Check filter indicators"
$ws2.Cells.Item(4, 5).Value = ""
$ws2.Cells.Item(4, 6).Value = "Two separate filter indicators are visible in the UI"
$ws2.Cells.Item(4, 7).Value = "6414a0cd67102fc717c034d7"
$ws2.Cells.Item(4, 8).Value = "This test case has been built by GenAI Workbench for XL import via Internal Importer."
$ws2.Cells.Item(4, 9).Value = "User is logged in as investigator and on transaction list page`nMultiple filters are successfully applied"
$ws2.Cells.Item(4, 10).Value = "Core"
$ws2.Cells.Item(4, 11).Value = "external"
$ws2.Cells.Item(4, 12).Value = "INVHUB-10821"
$ws2.Cells.Item(4, 13).Value = "blocks"
$ws2.Cells.Item(4, 14).Value = "GenAI_Test_Case"
$ws2.Cells.Item(4, 15).Value = "IM-5000"
$ws2.Cells.Item(4, 16).Value = "blocks"
$ws2.Cells.Item(4, 17).Value = "IM-3000"
$ws2.Cells.Item(4, 18).Value = "3 - Medium"
$ws2.Cells.Item(4, 19).Value = 37
$ws2.Cells.Item(4, 20).Value = "Release-1.0"
$ws2.Cells.Item(4, 21).Value = "Dublin"

# ---------------------------------------------------------------------------
# Sheet3: "Verify sort persistence with filtered results" -> "Verify filter
# persistence across navigation" (text-only changes, no new rows).
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("B2").Value = "Verify filter persistence across navigation"
$ws3.Range("D2").Value = "Navigate to a different page"
$ws3.Range("F2").Value = "Page navigation is successful"
$ws3.Range("I2").Value = "User has applied filters to transaction list`nFilter settings are preserved across navigation"

$ws3.Range("B3").Value = "Verify filter persistence across navigation"
$ws3.Range("D3").Value = "Return to transaction list page"
$ws3.Range("F3").Value = "Previously applied filters remain active and results are filtered accordingly"
$ws3.Range("I3").Value = "User has applied filters to transaction list`nFilter settings are preserved across navigation"

# ---------------------------------------------------------------------------
# Sheet4: "Verify multiple column filtering" -> "Verify empty results
# handling" (text-only changes, no new rows).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet4")

$ws4.Range("B2").Value = "Verify empty results handling"
$ws4.Range("D2").Value = "Apply filters that result in no matching transactions"
# A leading single quote is Excel's "treat as text" marker and gets consumed
# on assignment, so a doubled leading quote is used to keep one literal
# apostrophe at the start of the stored value.
$ws4.Range("F2").Value = "''No matching transactions' message is displayed"
$ws4.Range("I2").Value = "User is on transaction list page`nEmpty state handling is properly displayed"

$ws4.Range("B3").Value = "Verify empty results handling"
$ws4.Range("D3").Value = "Verify clear filters option"
$ws4.Range("F3").Value = "Option to clear filters is visible and clickable"
$ws4.Range("I3").Value = "User is on transaction list page`nEmpty state handling is properly displayed"

# ---------------------------------------------------------------------------
# Sheet5 is removed entirely.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Sheet5")
$null = $ws5.Delete()
